$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths ---
$ws.Columns.Item(4).ColumnWidth = 32.875
$ws.Columns.Item(10).ColumnWidth = 22
$ws.Columns.Item(11).ColumnWidth = 32.125
$ws.Columns.Item(13).ColumnWidth = 19.75

# --- Extend header styling (rows 1-2) into the two new trailing columns L, M ---
$ws.Range("L1:M2").Style = "Good"
$ws.Range("J2").Style = "Good"

# --- Row 3 styling: same "Good" look extended across the new columns, plus a medium bottom border for the whole header block ---
$ws.Range("B3:C3").Style = "Good"
$ws.Range("I3:M3").Style = "Good"
$r3 = $ws.Range("A3:M3")
$r3.Borders.Item(9).Weight = -4138

# --- New data row (Character / actor) ---
$ws.Range("C5").Value = "Character"
$ws.Range("B5").Value = "actor.TbCharacter"
$ws.Range("D5").Value = $true
$ws.Range("E5").Value = "Characters"
$ws.Range("I5").Value = "目录下包含所有英雄和敌人角色"

# --- Trailing blank formatted rows ---
$ws.Range("A6:M7").Style = "Normal"

$ws.Range("C10").Select()
